# Append two new result rows to the "Sonuçlar" sheet, mirroring the
# site's "personal download" export (rows 6 and 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 -----------------------------------------------------------
$ws.Range("A6").Value = "2025-08-11 16:49"
$ws.Range("B6").Value = "gh"
$ws.Range("C6").Value = "gh"
$ws.Range("D6").Value = 32
$ws.Range("E6").Value = 31
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = 27
$ws.Range("H6").Value = 31.5
$ws.Range("I6").Value = 28.5
$ws.Range("J6").Value = 0.3896
$ws.Range("K6").Value = 0.2666
$ws.Range("L6").Value = 0.2041
$ws.Range("M6").Value = 0.1396
$ws.Range("N6").Value = 38.96
$ws.Range("O6").Value = 26.66
$ws.Range("P6").Value = 20.41
$ws.Range("Q6").Value = 13.96

# --- Row 7 -------------------------------------------------------------
$ws.Range("A7").Value = "2025-08-11 17:06"
$ws.Range("B7").Value = "yaren"
$ws.Range("C7").Value = "cvf"
$ws.Range("D7").Value = 39
$ws.Range("E7").Value = 35
$ws.Range("F7").Value = 28
$ws.Range("G7").Value = 29
$ws.Range("H7").Value = 37
$ws.Range("I7").Value = 28.5
$ws.Range("J7").Value = 0.46
$ws.Range("K7").Value = 0.31
$ws.Range("L7").Value = 0.14
$ws.Range("M7").Value = 0.09
$ws.Range("N7").Value = "%45.77"
$ws.Range("O7").Value = "%31.32"
$ws.Range("P7").Value = "%13.61"
$ws.Range("Q7").Value = "%9.31"

Write-Host "Appended rows 6 and 7"
